$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (RF) - updated metrics
$ws.Range("B6").Value = 0.8219858156028369
$ws.Range("C6").Value = 0.8340575926717945
$ws.Range("D6").Value = 0.8219858156028369
$ws.Range("E6").Value = 0.824133966676615
$ws.Range("F6").Value = 0.8134294211850837
$ws.Range("G6").Value = 0.8211924573149492
$ws.Range("H6").Value = 0.8134294211850837
$ws.Range("I6").Value = 0.8138681890195073
$ws.Range("J6").Value = 0.7533287577213453
$ws.Range("K6").Value = 0.7644274484152398
$ws.Range("L6").Value = 0.7533287577213453
$ws.Range("M6").Value = 0.7534000713955357
$ws.Range("N6").Value = 0.8047357584076871
$ws.Range("O6").Value = 0.8137424017370091
$ws.Range("P6").Value = 0.8047357584076871
$ws.Range("Q6").Value = 0.8051615125966769
$ws.Range("R6").Value = 0.8197208876687258
$ws.Range("S6").Value = 0.8289062718262024
$ws.Range("T6").Value = 0.8197208876687258
$ws.Range("U6").Value = 0.8206418104533046
$ws.Range("V6").Value = 0.8133836650652025
$ws.Range("W6").Value = 0.8201578911488422
$ws.Range("X6").Value = 0.8133836650652025
$ws.Range("Y6").Value = 0.8125355882640468

# Row 7 (Ensemble) - updated metrics
$ws.Range("B7").Value = 0.8541066117593228
$ws.Range("C7").Value = 0.8576568759891664
$ws.Range("D7").Value = 0.8541066117593228
$ws.Range("E7").Value = 0.8542573678287834
$ws.Range("F7").Value = 0.8605811027224892
$ws.Range("G7").Value = 0.8639139817632244
$ws.Range("H7").Value = 0.8605811027224892
$ws.Range("I7").Value = 0.860939356328483
$ws.Range("J7").Value = 0.8262411347517731
$ws.Range("K7").Value = 0.8374450754791903
$ws.Range("L7").Value = 0.8262411347517731
$ws.Range("M7").Value = 0.8282070169752302
$ws.Range("N7").Value = 0.8776710134980554
$ws.Range("O7").Value = 0.8818057467490332
$ws.Range("P7").Value = 0.8776710134980554
$ws.Range("Q7").Value = 0.8773861715477101
$ws.Range("R7").Value = 0.8605811027224892
$ws.Range("S7").Value = 0.8661648814547099
$ws.Range("T7").Value = 0.8605811027224892
$ws.Range("U7").Value = 0.8599897504776187
$ws.Range("V7").Value = 0.8604438343628461
$ws.Range("W7").Value = 0.866328962456571
$ws.Range("X7").Value = 0.8604438343628461
$ws.Range("Y7").Value = 0.8606441137027421
